# Increase right margin in pipeline output files.
#
# In PowerPoint, shape position/size round-trips through single-precision
# (point-based) floats, so the EMU values stored in the OOXML after a resize
# can shift by a single EMU even on an axis that "shouldn't" have moved. To
# reproduce that faithfully we set BOTH .Left and .Width (in points) for each
# shape that was actually resized by the author, rather than only nudging the
# width. Shapes that are glued to these via connector site (stCxn/endCxn)
# pick up the new edge automatically, so we leave the connectors alone.

$EMU_PER_POINT = 12700

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "Output files" rounded rectangle (id=26) grows wider on the right.
$rr = $s.Shapes.Item("Rounded Rectangle 25")
$rr.Left = 6986723 / $EMU_PER_POINT
$rr.Width = 2188807 / $EMU_PER_POINT

# The two "Document" flow-chart shapes inside it (XGMML graph / id=145, and
# the one below it / id=146) also grow wider on the right.
$doc1 = $s.Shapes.Item("Document 144")
$doc1.Left = 7156170 / $EMU_PER_POINT
$doc1.Width = 1821979 / $EMU_PER_POINT

$doc2 = $s.Shapes.Item("Document 145")
$doc2.Left = 7156170 / $EMU_PER_POINT
$doc2.Width = 1821979 / $EMU_PER_POINT
